# Update fractions eating various staple foods to avoid 0 population
# in "Baseline year population inputs" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Baseline year population inputs")

# Set fraction eating rice (frac_rice) and fraction eating wheat (frac_wheat)
# from 0 to 0.1, letting the "other staples" formula (C19) recalculate
# automatically.
$ws.Range("C16").Value = 0.1
$ws.Range("C17").Value = 0.1

# Update the active selection to C18, matching the saved view state.
$ws.Activate()
$ws.Range("C18").Select()

$excel.CalculateFull()
